$wb = $excel.ActiveWorkbook

# --- Sheet "Feuil1" ---
$wsFeuil1 = $wb.Worksheets.Item("Feuil1")

# New base category row: RT.ART joins RO.ACT / RO.FOU / MP.CPT / ZZ.XXX.
$wsFeuil1.Range("D9").Value = "RT.ART"

# Update the remembered selection on this sheet.
$wsFeuil1.Range("D9").Select()

# --- Sheet "LIST" ---
$wsList = $wb.Worksheets.Item("LIST")

# A3 used to reference "RO.FOU.001.SUP.01"; it now becomes the new MAJ test case.
$wsList.Range("A3").Value = "RT.ART.001.MAJ"

# Rows 4 and 5 only held empty, formatted placeholder cells - clear them away entirely.
$wsList.Range("A4:A5").Clear()

# Update the remembered selection on this sheet.
$wsList.Range("B7:B8").Select()
$wsList.Application.ActiveCell.Worksheet.Range("B8").Activate()
